# Apply the edits described by the commit "Adjustments and Example Provided":
#   1. Refresh the cached "datetimeFigureOut" date placeholder text (6/8/2023 ->
#      3/24/2024) on the slide master and every slide layout.
#   2. Remove the "TextBox 13" shape from slide 1.

$p = $ppt.ActivePresentation

$newDate = "3/24/2024"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $isDatePlaceholder = $false
        if ($sh.Type -eq 14) {
            try {
                if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                    $isDatePlaceholder = $true
                }
            } catch {
                $isDatePlaceholder = $false
            }
        }
        if ($isDatePlaceholder) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# Slide master date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

# Every slide layout's date placeholder.
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}

# Remove the "TextBox 13" shape on slide 1.
$slide = $p.Slides.Item(1)
for ($i = $slide.Shapes.Count; $i -ge 1; $i--) {
    $sh = $slide.Shapes.Item($i)
    if ($sh.Name -eq "TextBox 13") {
        $sh.Delete()
    }
}
